$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix up existing rows 3 & 4 ---
# Row 3: rename task label (string content change only)
$ws.Range("B3").Value = "Package Rope Physics Exploration/Player Unit Tests"

# Row 4: rename task label and correct hours spent
$ws.Range("B4").Value = "Player Unit Tests"
$ws.Range("C4").Value = 4

# --- Append new log entries (rows 5-9) ---
# Copy the date style from A4 down to A5:A9 so new dates keep the same number format
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5:A9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A5").Value = 44259
$ws.Range("B5").Value = "Player Unit Tests"
$ws.Range("C5").Value = 4

$ws.Range("A6").Value = 44260
$ws.Range("B6").Value = "Dog Unit Tests"
$ws.Range("C6").Value = 2

$ws.Range("A7").Value = 44261
$ws.Range("B7").Value = "Dog Unit Tests"
$ws.Range("C7").Value = 4

$ws.Range("A8").Value = 44262
$ws.Range("B8").Value = "Dog Unit Tests"
$ws.Range("C8").Value = 4

$ws.Range("A9").Value = 44263
$ws.Range("B9").Value = "SpawnManager/Obstacle/GameManager Unit Tests"
$ws.Range("C9").Value = 4

# Move active selection to A10 like in the final saved state
$ws.Range("A10").Select() | Out-Null
